$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$texts = @(
  "Biz bu okula boşuna mı para veriyoruz veya verdiğimiz para nereye gidiyor, dk otomatları günlerdir bomboş, sebillerin çoğunda su bulunmuyor bunun sebebi nedir?",
  "Geçenlerde yurtta bir yemek yapalım dedik vaziyet bu, daha önce bir sürü onur kırıcı söz yazıldı, hakaretlerde bulunuldu, şikayet edildi ve halen yine aynı durumu yapmaya devam ediyorlar açık açık ne istediklerini söylesinler onlar da rahat etsin biz de rahat edelim.",
  "Çöpünüzü çöpe atmayı çoktan öğrenmiş olmanız gerekiyordu ama ne yazık ki öğrenememişsiniz. Kütüphanedeki camlı odaları lütfen pis bırakmayın kimse sizin çöpünüzü toplamak zorunda değil.",
  "Pandemi önlemi altında vize haftası yaptılar ancak bir sınıfta 60 kişi sınava giriyoruz. Herkes yanyana oturuyor ve minimum 2 saat aynı kalabalıkta öksüren hapşıran insan kalabalığı ile birlikteyiz. Önlemden kastınız bu mu?",
  "Birisi bana açıklayabilir mi aklım almıyor çünkü. Bu okul sanki dağ başında değilmiş gibi saat 8:30 a sınav koyuyorlar. Servise gidiyorsun serviste yer yok, otobüs seferleri kısıtlı. Lisede bile sınavlar en kötü ihtimal 2. derste yapılırdı. Merak ediyorum aklı başında bir yetkili bunları düşünemiyor mu? Kusura bakmayın ama bütün sorumluluk bizler de değil. Bu sınav programını yapanlar keşke biraz anlayışlı olsa.",
  "Sınav takvimini açıklamak yerine sınav haftası yurtlarda konaklama ücretini açıklayan bir okulumuz var. Tek dertleri para. Öğrenciler umurlarında değil maalesef.",
  "Merhabalar, acaba bu okul sınavlar ile ilgili herhangi bir açıklama yapacak mı? Derslerimin hepsi online ve kalacak yer bulamıyorum. Ben gibi bir çok arkadaşım da mağdur. Sınavlara 1 hafta kalmışken bize hiçbir açıklama yapılmıyor bu hafta haber vericez diyorlar neyi verecekler acaba? Son 3 gün kala bunca insan nerde kalacak, herhangi bir şeyden haberi olan varsa yazabilirse çok sevinirim. ",
  "Merhaba arkadaşlar bir şey sormak istiyorum. 1. sınıfta üstten ders almak için gerekli şartlar neler? Bilen varsa yardımcı olursa çok sevinirim teşekkürler.",
  "Merhaba, masa tenisi klübü hakkında bilgi almak istiyorum klüpte yetkili biri ulaşabilir mi ya da beğensin ben ulaşırım",
  "Merhabalar, okulu dondurmak istiyorum ne yapmam gerekiyor bilgi verebilir misiniz öğrenci işlerini aradım ama ulaşamadım bilgilendirirseniz sevinirim.",
  "Merhaba bana okuldan mail olarak 4333 TL tutarında e-arşiv faturası gönderilmiş. Mail olarak geri dönüş yaptım ancak cevap gelmedi. Bunun manası ne oluyor tam olarak? Başka bu maili alan var mı?",
  "Okulun shuttle saatileriyle ilgili kime gitmek gerekiyor. Cumartesi günü en geç 21:00'da maslaktan servis var, eskiden 00:00 servisleri olurdu. Bi rahatsız olan ben olamam bu durumdan. Konsey mi artık kim ilgileniyor bunula gidelim görüşelim.",
  "Arkadaşlar okulda ya da derslerde maske takmıyorsunuz, size yaklaşmıyoruz anladık ama servislerde takın bari. 1 saat yüzüne yapışmaz o maske merak etmeyin!",
  "Add drop haftası haftaya mı başlıyor. Staj raporlarını haftaya mı vermeye başlayacağız?",
  "Oryantasyon kapalı spor salonunda devam ediyo harika tedbir, hayet kalabalık ortamlar falan paylaşımlar yapılıyor. Festival havası var Işık'ta bizim mezuniyetimize niçin aynı özen gösterilmedi. ",
  "Ders programı hala belli olmayan var mı benim dışımda?",
  "Merhaba, pandemi şartlarındayız ve okul cidden anlamda çok pis. Cafeler, yemehaneler, sınıflar, lavaboları hiç demiyorum. Heryer çok pis ya öğrenciler pis ya da okul temizliğe dikkat etmiyor.",
  "Arkadaşlar son dersimi de bu yaz verdim, temmuzdan sonrasında işlendi sisteme. Öğrenci işleri diyor ki yapmanız gereken işlemler var vesaire. Tam ne diye soracakken hat kesiliyor, bir daha da bağlanamıyorum. Mezuniyet için yapmanız gereken işlemler neler? Bir de diyor ki seneye kaldı sizin mezuniyetiniz. Bilgi verebilecek var mı?",
  "Ders programımı yaptım ama bir sorundan dolayı onay için danışmana gönderemedim. Blackboard'umda dersler gözükmüyor ancak e-campus'te gözüküyor. Bana sorun çıkartır mı?",
  "Merhaba, birinci sınıfların paket program ataması hala yapılmadı mı? E-campus'e girildiği zaman ders programı gözükmüyor. Herkes de böyle mi?",
  "Okulun futbol takımı seçmeleri ne zaman oluyor bilen var mı acaba?",
  "OOC kredili ve kredisiz dersler ne oluyor tam olarak kısaca açıklar mısınız?",
  "İkinci doz aşımı yeni oldum en az 14 gün önceden yaptırmak gerekli miydi? Negatif pcr istenir mi okula girerken?",
  "Arkadaşlar Emine Ekin'den neden bin tane mail alıyorum? Bir tek bana mı böyle? Gerçekten mail kutum hep ondan gelen maillerle dolu. Bu diğer arkadaşlardan kaynaklanıyorsa lütfen dikkat edebilirler mi?",
  "Arkadaşlar yurt taksidi ödeme seçeneği nasıl oluyor? Yurt anlaşmalı bankalardan tahsil ediyor da biz bankaya mı aylık olarak ödüyoruz? ",
  "Merhabalar Işık'tan yatay geçişle başka bi okula geçip sonrasında tekrar Işık'a geçme gibi bir durum yaşayan var mı?",
  "Ohes4412 dersi örgün gözüküyor ama o ders normalde de online değil miydi? Ve bu derste çakışma sayılmıyordu diye hatırlıyorum doğru mu?",
  "Merhaba, Bilişsel Süreçler ya da Psikopataloji dersini üstten almak istiyorum. Kalsam aynı dersi ücret ödemeden kendi döneminde tekrar alabilir miyim?",
  "Ben hazırlık okuyordum kaldım ve okulumu  dondurup şubat ayındaki sınava kaydım donuk bir şekilde gireceğim fakat kayıt dondururken okula herhangi bir ücret ödeyecek miyim?",
  "Ders programında derslerim online görünüyor fakat excel tablosunda hibrit yazıyor hangisi doğru?"
)

$startRow = 50
$startId = 48
for ($i = 0; $i -lt $texts.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $startId + $i
  $ws.Cells.Item($row, 2).Value = $texts[$i]
}

$ws.Range("D84").Select()
